$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Underlying data change: ECs ligand (Plau) expression values and
# ECs receptor (Igf2r) expression values were updated with new TPM data.
# Rows 2-4 (A=ECs): G (avg expr), H (total expr)
$ws.Range("G2:G4").Value = 4.616372666666667
$ws.Range("H2:H4").Value = 13.849118

# Rows 2,5,8 (D=ECs): M (avg expr), N (total expr)
$ws.Range("M2").Value = 13.76827833333333
$ws.Range("N2").Value = 41.304835
$ws.Range("M5").Value = 13.76827833333333
$ws.Range("N5").Value = 41.304835
$ws.Range("M8").Value = 13.76827833333333
$ws.Range("N8").Value = 41.304835

# Now recompute derived columns I, J, O, P, Q, R, S, T for all data rows (2-10)
# based on the formulas:
#   I = G / sum(G for each unique sending cluster)
#   J = H / sum(H for each unique sending cluster)
#   O = M / sum(M for each unique target cluster)
#   P = N / sum(N for each unique target cluster)
#   Q = G * M ; R = H * N ; S = I * O ; T = J * P

# Sums across the three distinct sending clusters (rows 2,5,8 each represent one cluster)
$sumG = $ws.Range("G2").Value2 + $ws.Range("G5").Value2 + $ws.Range("G8").Value2
$sumH = $ws.Range("H2").Value2 + $ws.Range("H5").Value2 + $ws.Range("H8").Value2

# Sums across the three distinct target clusters (rows 2,3,4 each represent one cluster)
$sumM = $ws.Range("M2").Value2 + $ws.Range("M3").Value2 + $ws.Range("M4").Value2
$sumN = $ws.Range("N2").Value2 + $ws.Range("N3").Value2 + $ws.Range("N4").Value2

for ($r = 2; $r -le 10; $r++) {
    $g = $ws.Cells.Item($r, 7).Value2
    $h = $ws.Cells.Item($r, 8).Value2
    $m = $ws.Cells.Item($r, 13).Value2
    $n = $ws.Cells.Item($r, 14).Value2

    $i = $g / $sumG
    $j = $h / $sumH
    $o = $m / $sumM
    $p = $n / $sumN
    $q = $g * $m
    $rr = $h * $n
    $s = $i * $o
    $t = $j * $p

    $ws.Cells.Item($r, 9).Value = $i
    $ws.Cells.Item($r, 10).Value = $j
    $ws.Cells.Item($r, 15).Value = $o
    $ws.Cells.Item($r, 16).Value = $p
    $ws.Cells.Item($r, 17).Value = $q
    $ws.Cells.Item($r, 18).Value = $rr
    $ws.Cells.Item($r, 19).Value = $s
    $ws.Cells.Item($r, 20).Value = $t
}
